# ---------------------------------------------------------------------------
# 1. Re-style the table so it uses the built-in table style
#    {6CA22A20-C1E4-4463-AB03-EC456F8E0135} instead of the custom
#    "Table_0" style ({42BC6BAA-81F8-4F2E-A130-E037BA20A413}) that shipped
#    with the deck.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{6CA22A20-C1E4-4463-AB03-EC456F8E0135}", $false)
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Switch the deck's applied theme colours from the bespoke "Integral" /
#    "Red Violet" scheme to the stock "Office Theme" / "Office" scheme.
#    Order is dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (theme colour
#    indices 1-12). Values are VBA-style BGR-encoded "long" colours
#    (0xBBGGRR), i.e. the byte-reversal of the target RRGGBB hex values:
#      dk1      000000 -> 0x000000
#      lt1      FFFFFF -> 0xFFFFFF
#      dk2      44546A -> 0x6A5444
#      lt2      E7E6E6 -> 0xE6E6E7
#      accent1  5B9BD5 -> 0xD59B5B
#      accent2  ED7D31 -> 0x317DED
#      accent3  A5A5A5 -> 0xA5A5A5
#      accent4  FFC000 -> 0x00C0FF
#      accent5  4472C4 -> 0xC47244
#      accent6  70AD47 -> 0x47AD70
#      hlink    0563C1 -> 0xC16305
#      folHlink 954F72 -> 0x724F95
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,
    0xFFFFFF,
    0x6A5444,
    0xE6E6E7,
    0xD59B5B,
    0x317DED,
    0xA5A5A5,
    0x00C0FF,
    0xC47244,
    0x47AD70,
    0xC16305,
    0x724F95
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
